$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 22 only held a leftover empty/styled placeholder cell
# (C22); clear it completely so the new row can start from the column's
# default formatting, same as every other brand-new cell below.
$ws.Range("C22").Clear()

# New rows of translation data (id / Espanol / English), appended after the
# existing "essence" row (row 21): a prototype "damage" stat item and a
# prototype "hp to damage" stat item, plus the shared stat-effect-description
# strings used by their tooltips.

$ws.Range("A22").Value = "dmg_item_name"
$ws.Range("B22").Value = "Biceps"
$ws.Range("C22").Value = "Biceps"

$ws.Range("A23").Value = "dmg_item_desc"
$ws.Range("B23").Value = " "
$ws.Range("C23").Value = " "

$ws.Range("A24").Value = "stat_effect_desc"
$ws.Range("B24").Value = "Incrementa {0} en {1}{2}."
$ws.Range("C24").Value = "Increments {0} by {1}{2}."

$ws.Range("A25").Value = "refstat_effect_desc"
$ws.Range("B25").Value = "Incrementa tu {0} en {1}{2} de tu {3}."
$ws.Range("C25").Value = "Increments {0} by {1}{2} of your {3}."

$ws.Range("A26").Value = "hptodmg_item_name"
$ws.Range("B26").Value = "Gigantismo"
$ws.Range("C26").Value = "Giantism"

$ws.Range("A27").Value = "hptodmg_item_desc"
$ws.Range("B27").Value = " "
$ws.Range("C27").Value = " "

# Formatting touch-ups to match the sheet's existing conventions:
# - the "item name" id cells (same look as A14/A20's A cell) use style 3
# - the "desc" placeholder cells (same look as C15/C20) use style 4 on column C
$ws.Range("A14").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A26").PasteSpecial(-4122)

$ws.Range("C20").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D21").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A27").Select() | Out-Null
